$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.458.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.76%  "
$ws.Range("D3").Value = "'1.829.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'315.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").Value = "'0.5217"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.81%  "
$ws.Range("D8").Value = "'0.3911"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.51%  "
$ws.Range("D9").Value = "'0.07651"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.57%  "
$ws.Range("D10").Value = "'41.86"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.91%  "
$ws.Range("D11").Value = "'1.110"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.27%  "
$ws.Range("D12").Value = "'21.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.61%  "
$ws.Range("D13").Value = "'6.296"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.80%  "
$ws.Range("D14").Value = "'1.002"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").Value = "'7.549"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("D16").Value = "'1.824.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.80%  "
$ws.Range("D17").Value = "'93.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.76%  "
$ws.Range("D18").Value = "'0.00001083"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.21%  "
$ws.Range("D19").Value = "'0.06672"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.41%  "
$ws.Range("D20").Value = "'17.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.81%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").Value = "'6.189"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.83%  "
$ws.Range("D23").Value = "'28.475.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.73%  "
$ws.Range("D24").Value = "'11.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("E25").Value = "  +7.92%  "
$ws.Range("D26").Value = "'156.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  +2.27%  "
$ws.Range("D28").Value = "'2.035.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.75%  "
$ws.Range("D29").Value = "'2.398"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.02%  "
$ws.Range("D30").Value = "'125.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.58%  "
$ws.Range("D31").Value = "'1.122"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'0.1088"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.25%  "
$ws.Range("D33").Value = "'5.672"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.05%  "
$ws.Range("D34").Value = "'3.663"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.26%  "
$ws.Range("D35").Value = "'0.07016"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.79%  "
$ws.Range("D36").Value = "'0.2229"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.56%  "
$ws.Range("D37").Value = "'8.991"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.07%  "
$ws.Range("D38").Value = "'0.02325"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("D39").Value = "'5.141"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("D40").Value = "'0.6288"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.76%  "
$ws.Range("D41").Value = "'11.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "'1.183"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").Value = "'1.0000"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").Value = "'1.398"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.25%  "
$ws.Range("D45").Value = "'13.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("D46").Value = "'0.5905"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.24%  "
$ws.Range("D47").Value = "'3.714"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.12%  "
$ws.Range("D48").Value = "'124.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("D49").Value = "'1.981"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.16%  "
$ws.Range("E50").Value = "  +1.74%  "
$ws.Range("D51").Value = "'0.06933"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.91%  "
